$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B11: the Cedula value was stored as text; convert it to a real number ---
$ws.Range("B11").Value = 1000271912

# --- Append new log row 12 (registro actualizacion 2025-10-16 20:40:18 -> evento 2025-10-16 15:40:17) ---
$ws.Range("A12").Value = "2025-10-16 15:40:17"

# B12 keeps the Cedula as text (matches the source log format for this entry)
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "1000271912"

$ws.Range("C12").Value = "Maria"
$ws.Range("D12").Value = "TARJETA DE CRÉDITO"
$ws.Range("E12").Value = "****6898"
$ws.Range("F12").Value = "REESTRUCTURACION CON PAGO"
$ws.Range("G12").Value = "24 cuotas"
$ws.Range("H12").Value = "35.197.92.111"
$ws.Range("I12").Value = "The Dalles"
$ws.Range("J12").Value = "Oregon"
$ws.Range("K12").Value = "United States"
$ws.Range("L12").Value = "2025-10-16 15:40:17"
$ws.Range("M12").Value = "****6898"
$ws.Range("N12").Value = "35.197.92.111"
